$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 16.02.2022 09:30"

# Update row 8 (Benzina Albert Modřice): D8 delta and E8 last-checked timestamp
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 44608.38604166666
$ws.Range("E8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
